$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.652.35"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "3.399.26"
$ws.Range("E3").Value = "  -0.35%  "

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.08%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.59"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.53%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.42"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.61%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -2.88%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +0.02%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.719"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -1.55%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -6.91%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.35"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "3.942.57"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.07"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +1.68%  "

$ws.Range("E14").Value = "  -0.11%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000210"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -3.86%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.31"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("D17").Value = "3.409.37"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.23"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +1.52%  "

$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.07"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.68%  "

$ws.Range("D20").Value = "61.669.62"
$ws.Range("E20").Value = "  -0.42%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "473.29"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +16.38%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.87"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.39%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.24"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +2.11%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.98"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -0.73%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.27"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.93%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.60"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +10.73%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "32.78"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -0.29%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.76"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -0.72%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.83"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +3.26%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.83"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -0.27%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.62"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("E32").Value = "  -2.00%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -4.91%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.55"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -5.53%  "

$ws.Range("E35").Value = "  -0.82%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.78"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +4.94%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0483"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -3.02%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +0.19%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.02"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +5.15%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.321"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +3.25%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.133"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.16%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +8.13%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.04"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +3.94%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("E47").Value = "  +18.96%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.45"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -1.12%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.81"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("E50").Value = "  +7.36%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.56"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +12.93%  "
